$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI pipeline was rerun with new TPM data. The "Target cluster" = ECs
# rows are no longer produced, so drop them (original sheet rows 2, 5, 8, 11,
# i.e. one row per "Sending cluster" group). Delete bottom-up so the row
# indices of the rows still to be removed stay valid.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# The remaining 8 rows keep their relative order but get recomputed
# expression/specificity statistics from the new TPM run. Write them in one
# shot as a 2D array (rows x columns A:T).
$data = New-Object 'object[,]' 8,20
$data[0,0] = "ECs"
$data[0,1] = "Efna2"
$data[0,2] = "Epha5"
$data[0,3] = "FAPs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 3.099562666666667
$data[0,7] = 9.298688
$data[0,8] = 0.2686390288432488
$data[0,9] = 0.2686390288432488
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.06103333333333334
$data[0,13] = 0.1831
$data[0,14] = 0.2094151016766933
$data[0,15] = 0.2094151016766933
$data[0,16] = 0.1891766414222222
$data[0,17] = 1.7025897728
$data[0,18] = 0.05625706953953707
$data[0,19] = 0.05625706953953707
$data[1,0] = "ECs"
$data[1,1] = "Efna2"
$data[1,2] = "Epha5"
$data[1,3] = "MuSCs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 3.099562666666667
$data[1,7] = 9.298688
$data[1,8] = 0.2686390288432488
$data[1,9] = 0.2686390288432488
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.2304133333333333
$data[1,13] = 0.69124
$data[1,14] = 0.7905848983233067
$data[1,15] = 0.7905848983233067
$data[1,16] = 0.7141805659022222
$data[1,17] = 6.42762509312
$data[1,18] = 0.2123819593037117
$data[1,19] = 0.2123819593037117
$data[2,0] = "FAPs"
$data[2,1] = "Efna2"
$data[2,2] = "Epha5"
$data[2,3] = "FAPs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 6.189892666666666
$data[2,7] = 18.569678
$data[2,8] = 0.5364778626674904
$data[2,9] = 0.5364778626674905
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.06103333333333334
$data[2,13] = 0.1831
$data[2,14] = 0.2094151016766933
$data[2,15] = 0.2094151016766933
$data[2,16] = 0.3777897824222222
$data[2,17] = 3.4001080418
$data[2,18] = 0.1123465661578076
$data[2,19] = 0.1123465661578076
$data[3,0] = "FAPs"
$data[3,1] = "Efna2"
$data[3,2] = "Epha5"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 6.189892666666666
$data[3,7] = 18.569678
$data[3,8] = 0.5364778626674904
$data[3,9] = 0.5364778626674905
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.2304133333333333
$data[3,13] = 0.69124
$data[3,14] = 0.7905848983233067
$data[3,15] = 0.7905848983233067
$data[3,16] = 1.426233802302222
$data[3,17] = 12.83610422072
$data[3,18] = 0.4241312965096828
$data[3,19] = 0.4241312965096829
$data[4,0] = "MuSCs"
$data[4,1] = "Efna2"
$data[4,2] = "Epha5"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.888584
$data[4,7] = 5.665752
$data[4,8] = 0.1636835341659699
$data[4,9] = 0.1636835341659699
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.06103333333333334
$data[4,13] = 0.1831
$data[4,14] = 0.2094151016766933
$data[4,15] = 0.2094151016766933
$data[4,16] = 0.1152665768
$data[4,17] = 1.0373991912
$data[4,18] = 0.03427780395016709
$data[4,19] = 0.03427780395016709
$data[5,0] = "MuSCs"
$data[5,1] = "Efna2"
$data[5,2] = "Epha5"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.888584
$data[5,7] = 5.665752
$data[5,8] = 0.1636835341659699
$data[5,9] = 0.1636835341659699
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.2304133333333333
$data[5,13] = 0.69124
$data[5,14] = 0.7905848983233067
$data[5,15] = 0.7905848983233067
$data[5,16] = 0.43515493472
$data[5,17] = 3.91639441248
$data[5,18] = 0.1294057302158028
$data[5,19] = 0.1294057302158028
$data[6,0] = "Resolving-Mac"
$data[6,1] = "Efna2"
$data[6,2] = "Epha5"
$data[6,3] = "FAPs"
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.3599813333333333
$data[6,7] = 1.079944
$data[6,8] = 0.03119957432329092
$data[6,9] = 0.03119957432329093
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.06103333333333334
$data[6,13] = 0.1831
$data[6,14] = 0.2094151016766933
$data[6,15] = 0.2094151016766933
$data[6,16] = 0.02197086071111111
$data[6,17] = 0.1977377464
$data[6,18] = 0.006533662029181517
$data[6,19] = 0.006533662029181518
$data[7,0] = "Resolving-Mac"
$data[7,1] = "Efna2"
$data[7,2] = "Epha5"
$data[7,3] = "MuSCs"
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 0.3599813333333333
$data[7,7] = 1.079944
$data[7,8] = 0.03119957432329092
$data[7,9] = 0.03119957432329093
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.2304133333333333
$data[7,13] = 0.69124
$data[7,14] = 0.7905848983233067
$data[7,15] = 0.7905848983233067
$data[7,16] = 0.08294449895111111
$data[7,17] = 0.74650049056
$data[7,18] = 0.02466591229410941
$data[7,19] = 0.02466591229410941

$ws.Range("A2:T9").Value2 = $data
